# Rebrand the "Product Development" KPI dashboard template to an
# "Artificial Intelligence and Machine Learning" KPI dashboard template.

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions & User Guide")
$wsDashboard    = $wb.Worksheets.Item("KPI Dashboard")

# ---- Instructions & User Guide sheet ----------------------------------
$wsInstructions.Range("A1").Value = "Artificial Intelligence and Machine Learning KPI Dashboard - User Guide & Instructions"
$wsInstructions.Range("B23").Value = "Availability and reliability of AI systems"

# ---- KPI Dashboard sheet ------------------------------------------------
$wsDashboard.Range("A1").Value = "ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING - KPI DASHBOARD"
$wsDashboard.Range("A2").Value = "Project: AI/ML Implementation"

# Update the recurring "Notes" column (K) for every KPI row (8-22)
$kpiNote = "Critical KPI for Artificial Intelligence and Machine Learning success"
for ($row = 8; $row -le 22; $row++) {
    $wsDashboard.Cells.Item($row, 11).Value = $kpiNote
}

# Update "Owner" column (I) rows that referenced Product Engineers
$wsDashboard.Range("I10").Value = "ML Engineers"
$wsDashboard.Range("I16").Value = "ML Engineers"
$wsDashboard.Range("I22").Value = "ML Engineers"
